$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Swap "Groenlandia" and "Islas Malvinas" rows (A210 / A211)
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"

# Update "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 07:43"

# Row 15 - Pakistan
$ws.Range("B15").Value = 271887
$ws.Range("C15").Value = 1487
$ws.Range("D15").Value = 236596
$ws.Range("E15").Value = 29504
$ws.Range("G15").Value = 24
$ws.Range("H15").Value = 5787

# Row 65 - Uzbekistan
$ws.Range("B65").Value = 19653
$ws.Range("C65").Value = 293
$ws.Range("E65").Value = 9074
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = 107

# Row 75 - El Salvador
$ws.Range("D75").Value = 7455
$ws.Range("E75").Value = 5947
$ws.Range("G75").Value = 11
$ws.Range("H75").Value = 390

# Row 108 - Tailandia
$ws.Range("B108").Value = 3282
$ws.Range("C108").Value = 3
$ws.Range("D108").Value = 3109
$ws.Range("E108").Value = 115

# Row 187 - Butan
$ws.Range("D187").Value = 85
$ws.Range("E187").Value = 7
